# ----------------------------------------------------------------------------
# Rename the "Requested quantity" headers so they describe each report more
# specifically, and add a brand-new "PO Forecast" sheet with the forecast
# data (ds / PO_Forecast / yhat_lower / yhat_upper).
# ----------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$ws2 = $wb.Worksheets.Item(2)   # "Monthly Trend"

$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet right after "Monthly Trend" so the final sheet order is
# Weekly Quantity, Monthly Trend, PO Forecast.
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# Match page margins used by the other two sheets (0.75"/0.75"/1"/1", 0.5"/0.5").
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# Header row
$ws3.Cells.Item(1,1).Value = "ds"
$ws3.Cells.Item(1,2).Value = "PO_Forecast"
$ws3.Cells.Item(1,3).Value = "yhat_lower"
$ws3.Cells.Item(1,4).Value = "yhat_upper"

# Forecast data rows (46 rows, columns ds / PO_Forecast / yhat_lower / yhat_upper)
$ws3.Cells.Item(2,1).Value = 45172.99999999999
$ws3.Cells.Item(2,2).Value = 110
$ws3.Cells.Item(2,3).Value = -87.00597192098942
$ws3.Cells.Item(2,4).Value = 336.3399641993224
$ws3.Cells.Item(3,1).Value = 45179.99999999999
$ws3.Cells.Item(3,2).Value = 112
$ws3.Cells.Item(3,3).Value = -85.30664452171061
$ws3.Cells.Item(3,4).Value = 316.8940940530652
$ws3.Cells.Item(4,1).Value = 45186.99999999999
$ws3.Cells.Item(4,2).Value = 114
$ws3.Cells.Item(4,3).Value = -102.6498713576774
$ws3.Cells.Item(4,4).Value = 304.4735730352526
$ws3.Cells.Item(5,1).Value = 45200.99999999999
$ws3.Cells.Item(5,2).Value = 118
$ws3.Cells.Item(5,3).Value = -83.87881389937914
$ws3.Cells.Item(5,4).Value = 314.1242682937138
$ws3.Cells.Item(6,1).Value = 45207.99999999999
$ws3.Cells.Item(6,2).Value = 119
$ws3.Cells.Item(6,3).Value = -75.32888129781993
$ws3.Cells.Item(6,4).Value = 334.2648668910668
$ws3.Cells.Item(7,1).Value = 45214.99999999999
$ws3.Cells.Item(7,2).Value = 121
$ws3.Cells.Item(7,3).Value = -96.82167354694522
$ws3.Cells.Item(7,4).Value = 327.9605789595081
$ws3.Cells.Item(8,1).Value = 45221.99999999999
$ws3.Cells.Item(8,2).Value = 123
$ws3.Cells.Item(8,3).Value = -84.72169157018998
$ws3.Cells.Item(8,4).Value = 315.8010178017082
$ws3.Cells.Item(9,1).Value = 45228.99999999999
$ws3.Cells.Item(9,2).Value = 125
$ws3.Cells.Item(9,3).Value = -66.78636098768818
$ws3.Cells.Item(9,4).Value = 337.6486761917193
$ws3.Cells.Item(10,1).Value = 45242.99999999999
$ws3.Cells.Item(10,2).Value = 129
$ws3.Cells.Item(10,3).Value = -68.35876519604658
$ws3.Cells.Item(10,4).Value = 333.1527994082277
$ws3.Cells.Item(11,1).Value = 45249.99999999999
$ws3.Cells.Item(11,2).Value = 131
$ws3.Cells.Item(11,3).Value = -70.431826289135
$ws3.Cells.Item(11,4).Value = 349.7222557071648
$ws3.Cells.Item(12,1).Value = 45256.99999999999
$ws3.Cells.Item(12,2).Value = 133
$ws3.Cells.Item(12,3).Value = -82.29160817212585
$ws3.Cells.Item(12,4).Value = 328.0611031545547
$ws3.Cells.Item(13,1).Value = 45263.99999999999
$ws3.Cells.Item(13,2).Value = 135
$ws3.Cells.Item(13,3).Value = -65.9664887213377
$ws3.Cells.Item(13,4).Value = 346.4153779000004
$ws3.Cells.Item(14,1).Value = 45298.99999999999
$ws3.Cells.Item(14,2).Value = 145
$ws3.Cells.Item(14,3).Value = -60.73439032213153
$ws3.Cells.Item(14,4).Value = 354.8442395410246
$ws3.Cells.Item(15,1).Value = 45312.99999999999
$ws3.Cells.Item(15,2).Value = 148
$ws3.Cells.Item(15,3).Value = -62.4542711586262
$ws3.Cells.Item(15,4).Value = 349.6946127088023
$ws3.Cells.Item(16,1).Value = 45333.99999999999
$ws3.Cells.Item(16,2).Value = 154
$ws3.Cells.Item(16,3).Value = -50.77348457682366
$ws3.Cells.Item(16,4).Value = 367.1299307370362
$ws3.Cells.Item(17,1).Value = 45340.99999999999
$ws3.Cells.Item(17,2).Value = 156
$ws3.Cells.Item(17,3).Value = -43.08534283399006
$ws3.Cells.Item(17,4).Value = 368.679460430803
$ws3.Cells.Item(18,1).Value = 45347.99999999999
$ws3.Cells.Item(18,2).Value = 158
$ws3.Cells.Item(18,3).Value = -41.75375148854484
$ws3.Cells.Item(18,4).Value = 365.203815069939
$ws3.Cells.Item(19,1).Value = 45354.99999999999
$ws3.Cells.Item(19,2).Value = 160
$ws3.Cells.Item(19,3).Value = -41.47148752118873
$ws3.Cells.Item(19,4).Value = 376.6894389833689
$ws3.Cells.Item(20,1).Value = 45375.99999999999
$ws3.Cells.Item(20,2).Value = 166
$ws3.Cells.Item(20,3).Value = -35.77306484644606
$ws3.Cells.Item(20,4).Value = 372.7805506937665
$ws3.Cells.Item(21,1).Value = 45396.99999999999
$ws3.Cells.Item(21,2).Value = 172
$ws3.Cells.Item(21,3).Value = -41.40187525682166
$ws3.Cells.Item(21,4).Value = 379.7138185683627
$ws3.Cells.Item(22,1).Value = 45403.99999999999
$ws3.Cells.Item(22,2).Value = 174
$ws3.Cells.Item(22,3).Value = -37.48409680158849
$ws3.Cells.Item(22,4).Value = 370.6771232189626
$ws3.Cells.Item(23,1).Value = 45410.99999999999
$ws3.Cells.Item(23,2).Value = 175
$ws3.Cells.Item(23,3).Value = -33.52678067537869
$ws3.Cells.Item(23,4).Value = 386.0167175826879
$ws3.Cells.Item(24,1).Value = 45417.99999999999
$ws3.Cells.Item(24,2).Value = 177
$ws3.Cells.Item(24,3).Value = -40.22724801639843
$ws3.Cells.Item(24,4).Value = 389.2227467309363
$ws3.Cells.Item(25,1).Value = 45424.99999999999
$ws3.Cells.Item(25,2).Value = 179
$ws3.Cells.Item(25,3).Value = -15.89653875288731
$ws3.Cells.Item(25,4).Value = 390.3284586593616
$ws3.Cells.Item(26,1).Value = 45431.99999999999
$ws3.Cells.Item(26,2).Value = 181
$ws3.Cells.Item(26,3).Value = -15.8114382455161
$ws3.Cells.Item(26,4).Value = 386.9331808825547
$ws3.Cells.Item(27,1).Value = 45445.99999999999
$ws3.Cells.Item(27,2).Value = 185
$ws3.Cells.Item(27,3).Value = -14.19403495873571
$ws3.Cells.Item(27,4).Value = 383.4323209786514
$ws3.Cells.Item(28,1).Value = 45452.99999999999
$ws3.Cells.Item(28,2).Value = 187
$ws3.Cells.Item(28,3).Value = -31.19570163551963
$ws3.Cells.Item(28,4).Value = 385.5816976350756
$ws3.Cells.Item(29,1).Value = 45501.99999999999
$ws3.Cells.Item(29,2).Value = 201
$ws3.Cells.Item(29,3).Value = 6.110440775098745
$ws3.Cells.Item(29,4).Value = 408.2802252504969
$ws3.Cells.Item(30,1).Value = 45508.99999999999
$ws3.Cells.Item(30,2).Value = 203
$ws3.Cells.Item(30,3).Value = 12.29176367684655
$ws3.Cells.Item(30,4).Value = 399.7282938012293
$ws3.Cells.Item(31,1).Value = 45515.99999999999
$ws3.Cells.Item(31,2).Value = 204
$ws3.Cells.Item(31,3).Value = 11.2080838655321
$ws3.Cells.Item(31,4).Value = 427.8943998963074
$ws3.Cells.Item(32,1).Value = 45522.99999999999
$ws3.Cells.Item(32,2).Value = 206
$ws3.Cells.Item(32,3).Value = 5.872529928779913
$ws3.Cells.Item(32,4).Value = 420.1494496798093
$ws3.Cells.Item(33,1).Value = 45529.99999999999
$ws3.Cells.Item(33,2).Value = 208
$ws3.Cells.Item(33,3).Value = 14.82391323705306
$ws3.Cells.Item(33,4).Value = 420.2396645204598
$ws3.Cells.Item(34,1).Value = 45557.99999999999
$ws3.Cells.Item(34,2).Value = 216
$ws3.Cells.Item(34,3).Value = 7.933591528295727
$ws3.Cells.Item(34,4).Value = 411.3760033277171
$ws3.Cells.Item(35,1).Value = 45578.99999999999
$ws3.Cells.Item(35,2).Value = 222
$ws3.Cells.Item(35,3).Value = 12.60622419592377
$ws3.Cells.Item(35,4).Value = 436.8882725338094
$ws3.Cells.Item(36,1).Value = 45585.99999999999
$ws3.Cells.Item(36,2).Value = 224
$ws3.Cells.Item(36,3).Value = 13.90292832834669
$ws3.Cells.Item(36,4).Value = 430.2635849793872
$ws3.Cells.Item(37,1).Value = 45592.99999999999
$ws3.Cells.Item(37,2).Value = 226
$ws3.Cells.Item(37,3).Value = 10.34800177920913
$ws3.Cells.Item(37,4).Value = 446.5690149659062
$ws3.Cells.Item(38,1).Value = 45599.99999999999
$ws3.Cells.Item(38,2).Value = 228
$ws3.Cells.Item(38,3).Value = 13.80651517690695
$ws3.Cells.Item(38,4).Value = 435.0111690111014
$ws3.Cells.Item(39,1).Value = 45641.99999999999
$ws3.Cells.Item(39,2).Value = 239
$ws3.Cells.Item(39,3).Value = 33.98525882652864
$ws3.Cells.Item(39,4).Value = 438.6562434605662
$ws3.Cells.Item(40,1).Value = 45648.99999999999
$ws3.Cells.Item(40,2).Value = 241
$ws3.Cells.Item(40,3).Value = 35.73218350605679
$ws3.Cells.Item(40,4).Value = 448.4919077354782
$ws3.Cells.Item(41,1).Value = 45655.99999999999
$ws3.Cells.Item(41,2).Value = 243
$ws3.Cells.Item(41,3).Value = 42.95045864317697
$ws3.Cells.Item(41,4).Value = 451.0183937701397
$ws3.Cells.Item(42,1).Value = 45662.99999999999
$ws3.Cells.Item(42,2).Value = 245
$ws3.Cells.Item(42,3).Value = 39.97670994319832
$ws3.Cells.Item(42,4).Value = 448.0998433327757
$ws3.Cells.Item(43,1).Value = 45669.99999999999
$ws3.Cells.Item(43,2).Value = 247
$ws3.Cells.Item(43,3).Value = 62.91546732639416
$ws3.Cells.Item(43,4).Value = 447.0041701985634
$ws3.Cells.Item(44,1).Value = 45676.99999999999
$ws3.Cells.Item(44,2).Value = 249
$ws3.Cells.Item(44,3).Value = 52.44084654148256
$ws3.Cells.Item(44,4).Value = 464.5502637848652
$ws3.Cells.Item(45,1).Value = 45683.99999999999
$ws3.Cells.Item(45,2).Value = 251
$ws3.Cells.Item(45,3).Value = 35.96527835852763
$ws3.Cells.Item(45,4).Value = 459.7511116693439
$ws3.Cells.Item(46,1).Value = 45690.99999999999
$ws3.Cells.Item(46,2).Value = 253
$ws3.Cells.Item(46,3).Value = 52.36115022904327
$ws3.Cells.Item(46,4).Value = 462.6180795606361
$ws3.Cells.Item(47,1).Value = 45697.99999999999
$ws3.Cells.Item(47,2).Value = 255
$ws3.Cells.Item(47,3).Value = 46.20842101402304
$ws3.Cells.Item(47,4).Value = 461.0817698926137

# Re-use the same cell styles as the other sheets: bold/bordered/centered
# header row, and the date number format on column A.
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws1.Range("A2").Copy()
$ws3.Range("A2:A47").PasteSpecial(-4122)

# Keep the originally active sheet selected.
$ws1.Activate()
